# Update "想去人数" (F column) figures across the workbook's sheets to
# reflect newly generated output data (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1178
$ws1.Range("F4").Value  = 12606
$ws1.Range("F5").Value  = 708
$ws1.Range("F10").Value = 321
$ws1.Range("F16").Value = 120
$ws1.Range("F17").Value = 331
$ws1.Range("F19").Value = 286
$ws1.Range("F20").Value = 118
$ws1.Range("F23").Value = 207
$ws1.Range("F25").Value = 1255

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 275
$ws2.Range("F5").Value  = 4448
$ws2.Range("F6").Value  = 130
$ws2.Range("F10").Value = 350

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value  = 861

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 861
$ws4.Range("F6").Value  = 1178
$ws4.Range("F7").Value  = 12606
$ws4.Range("F8").Value  = 275
$ws4.Range("F9").Value  = 708
$ws4.Range("F14").Value = 321
$ws4.Range("F19").Value = 4448
$ws4.Range("F21").Value = 130
$ws4.Range("F22").Value = 130
$ws4.Range("F24").Value = 120
$ws4.Range("F28").Value = 350
$ws4.Range("F29").Value = 331
$ws4.Range("F32").Value = 286
$ws4.Range("F33").Value = 118
$ws4.Range("F37").Value = 207
$ws4.Range("F41").Value = 1255

$wb.Save()
